$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 568.5417
$ws.Range("I2").Value = 209.5
$ws.Range("J2").Value = 1286.625
$ws.Range("K2").Value = 209.5
$ws.Range("L2").Value = 1286.625
$ws.Range("M2").Value = -96.5
$ws.Range("N2").Value = -1512.625

$ws.Range("H5").Value = 544.0833
$ws.Range("I5").Value = 92.666664
$ws.Range("J5").Value = 1898.3334
$ws.Range("K5").Value = 92.666664
$ws.Range("L5").Value = 1898.3334
$ws.Range("M5").Value = 22.333336
$ws.Range("N5").Value = -2128.3334

$ws.Range("H8").Value = 754
$ws.Range("I8").Value = 22
$ws.Range("J8").Value = 2950
$ws.Range("K8").Value = 66
$ws.Range("L8").Value = 8850
$ws.Range("M8").Value = 73
$ws.Range("N8").Value = -9128

$ws.Range("H18").Value = 2116
$ws.Range("I18").Value = 1939.4
$ws.Range("J18").Value = 2999
$ws.Range("K18").Value = 1939.4
$ws.Range("L18").Value = 2999
$ws.Range("M18").Value = -1655.4
$ws.Range("N18").Value = -3567

$ws.Range("H32").Value = 2970.3635
$ws.Range("I32").Value = 1738.2
$ws.Range("J32").Value = 3997.1667
$ws.Range("K32").Value = 1738.2
$ws.Range("L32").Value = 3997.1667
$ws.Range("M32").Value = -1412.2
$ws.Range("N32").Value = -4649.1667

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws.Range("H69").Value = 30922.37
$ws.Range("I69").Value = 53102
$ws.Range("J69").Value = 23001.072
$ws.Range("K69").Value = 159306
$ws.Range("L69").Value = 69003.216
$ws.Range("M69").Value = -158432
$ws.Range("N69").Value = -70751.216

$ws.Range("H72").Value = 30922.37
$ws.Range("I72").Value = 53102
$ws.Range("J72").Value = 23001.072
$ws.Range("K72").Value = 477918
$ws.Range("L72").Value = 207009.648
$ws.Range("M72").Value = -473550
$ws.Range("N72").Value = -215745.648

$ws.Range("H116").Value = 2449.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2449.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 2449.5
$ws.Range("N116").Value = -9333.5

$ws.Range("H135").Value = 575.75
$ws.Range("I135").Value = 575.75
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5181.75
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2646.75

$ws.Range("H141").Value = 3249.5
$ws.Range("I141").Value = 2332.6667
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 6998.000100000001
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = -1818.000100000001
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2787.4
$ws.Range("I61").Value = 2919.4443
$ws.Range("J61").Value = 1599
$ws.Range("K61").Value = 2919.4443
$ws.Range("L61").Value = 1599
$ws.Range("M61").Value = -2707.4443
$ws.Range("N61").Value = -2023

$ws.Range("H102").Value = 952.5
$ws.Range("I102").Value = 952.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 952.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 669.5

$ws.Range("H136").Value = 2787.4
$ws.Range("I136").Value = 2919.4443
$ws.Range("J136").Value = 1599
$ws.Range("K136").Value = 8758.332900000001
$ws.Range("L136").Value = 4797
$ws.Range("M136").Value = -6208.332900000001
$ws.Range("N136").Value = -9897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1447975
$ws.Range("I7").Value = 1246.5
$ws.Range("J7").Value = 2026666.4
$ws.Range("K7").Value = 1246.5
$ws.Range("L7").Value = 2026666.4
$ws.Range("M7").Value = -1133.5
$ws.Range("N7").Value = -2026892.4

$ws.Range("H86").Value = 2240
$ws.Range("I86").Value = 2266.6667
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2266.6667
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -1143.6667
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 2240
$ws.Range("I89").Value = 2266.6667
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 11333.3335
$ws.Range("L89").Value = 2000
$ws.Range("M89").Value = -5717.333500000001
$ws.Range("N89").Value = -21232

$ws.Range("H99").Value = 1618.25
$ws.Range("I99").Value = 1491
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1491
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = 7
$ws.Range("N99").Value = -4996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2624
$ws.Range("I16").Value = 3036.5
$ws.Range("J16").Value = 1964
$ws.Range("K16").Value = 3036.5
$ws.Range("L16").Value = 1964
$ws.Range("M16").Value = -2749.5
$ws.Range("N16").Value = -2538

$ws.Range("H31").Value = 1507.5385
$ws.Range("I31").Value = 1074.875
$ws.Range("J31").Value = 2199.8
$ws.Range("K31").Value = 1074.875
$ws.Range("L31").Value = 2199.8
$ws.Range("M31").Value = -779.875
$ws.Range("N31").Value = -2789.8

$ws.Range("H34").Value = 1507.5385
$ws.Range("I34").Value = 1074.875
$ws.Range("J34").Value = 2199.8
$ws.Range("K34").Value = 1074.875
$ws.Range("L34").Value = 2199.8
$ws.Range("M34").Value = -872.875
$ws.Range("N34").Value = -2603.8

$ws.Range("H113").Value = 2624
$ws.Range("I113").Value = 3036.5
$ws.Range("J113").Value = 1964
$ws.Range("K113").Value = 3036.5
$ws.Range("L113").Value = 1964
$ws.Range("M113").Value = -866.5
$ws.Range("N113").Value = -6304

$ws.Range("H132").Value = 1880.3182
$ws.Range("I132").Value = 1816.2222
$ws.Range("J132").Value = 2168.75
$ws.Range("K132").Value = 5448.6666
$ws.Range("L132").Value = 6506.25
$ws.Range("M132").Value = -2918.6666
$ws.Range("N132").Value = -11566.25

$ws.Range("H141").Value = 165497.12
$ws.Range("I141").Value = 44989
$ws.Range("J141").Value = 182712.58
$ws.Range("K141").Value = 44989
$ws.Range("L141").Value = 182712.58
$ws.Range("M141").Value = -39809
$ws.Range("N141").Value = -193072.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H103").Value = 797.7778
$ws.Range("I103").Value = 36.8
$ws.Range("J103").Value = 1749
$ws.Range("K103").Value = 110.4
$ws.Range("L103").Value = 5247
$ws.Range("M103").Value = 768.6
$ws.Range("N103").Value = -7005

$ws.Range("H141").Value = 7453.222
$ws.Range("I141").Value = 7384.875
$ws.Range("J141").Value = 8000
$ws.Range("K141").Value = 22154.625
$ws.Range("L141").Value = 24000
$ws.Range("M141").Value = -16974.625
$ws.Range("N141").Value = -34360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.28125
$ws.Range("I2").Value = 88.333336
$ws.Range("J2").Value = 17.6
$ws.Range("K2").Value = 88.333336
$ws.Range("L2").Value = 17.6
$ws.Range("M2").Value = 24.666664
$ws.Range("N2").Value = -243.6

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H122").Value = 2380.077
$ws.Range("I122").Value = 2236.625
$ws.Range("J122").Value = 2609.6
$ws.Range("K122").Value = 6709.875
$ws.Range("L122").Value = 7828.799999999999
$ws.Range("M122").Value = -4259.875
$ws.Range("N122").Value = -12728.8

$ws.Range("H140").Value = 143843.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 143843.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 143843.5
$ws.Range("N140").Value = -154203.5

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1228.5714
$ws.Range("I22").Value = 1033.3334
$ws.Range("J22").Value = 1375
$ws.Range("K22").Value = 1033.3334
$ws.Range("L22").Value = 1375
$ws.Range("M22").Value = -738.3334
$ws.Range("N22").Value = -1965

$ws.Range("H27").Value = 1228.5714
$ws.Range("I27").Value = 1033.3334
$ws.Range("J27").Value = 1375
$ws.Range("K27").Value = 1033.3334
$ws.Range("L27").Value = 1375
$ws.Range("M27").Value = -926.3334
$ws.Range("N27").Value = -1589

$ws.Range("H46").Value = 3751.6667
$ws.Range("I46").Value = 2270.2
$ws.Range("J46").Value = 4492.4
$ws.Range("K46").Value = 2270.2
$ws.Range("L46").Value = 4492.4
$ws.Range("M46").Value = -2082.2
$ws.Range("N46").Value = -4868.4

$ws.Range("H55").Value = 574.4286
$ws.Range("I55").Value = 586.8333
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 586.8333
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -413.8333
$ws.Range("N55").Value = -846

$ws.Range("H136").Value = 1787.5
$ws.Range("I136").Value = 1787.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5362.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2812.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19993.5
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 19993.5
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 19993.5
$ws.Range("N41").Value = -20773.5

$ws.Range("H46").Value = 45412.223
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 45412.223
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 45412.223
$ws.Range("N46").Value = -45874.223

$ws.Range("H96").Value = 1700
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1700
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1700
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4446

$ws.Range("H100").Value = 20001094
$ws.Range("I100").Value = 25001192
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 50002384
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -50001843
$ws.Range("N100").Value = -2482

$ws.Range("H101").Value = 11959
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 11959
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 11959
$ws.Range("N101").Value = -18449

$ws.Range("H132").Value = 1935.4
$ws.Range("I132").Value = 2425.1667
$ws.Range("J132").Value = 1200.75
$ws.Range("K132").Value = 7275.500100000001
$ws.Range("L132").Value = 3602.25
$ws.Range("M132").Value = -4745.500100000001
$ws.Range("N132").Value = -8662.25

$ws.Range("H134").Value = 45412.223
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 45412.223
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 136236.669
$ws.Range("N134").Value = -141306.669
